$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 13 ("Are we ready?"),
# shifting "Are we ready?" and "????" down to rows 15-16.
$ws.Rows.Item(13).Insert()
$ws.Rows.Item(13).Insert()

$ws.Range("A13").Value = "12. what I learnt from Twitter reliability to my softwares..."
$ws.Range("A14").Value = "13. If I make a social media site what should I do to protect the reliability?"

$ws.Range("A21").Select()
